$d = $word.ActiveDocument

# The highlighted command line currently reads:
#   " artisan serve --host=0.0.0.0 --port="
# It needs to become " artisan serve " + "--port=" (i.e. remove the
# "--host=0.0.0.0 " portion), with the hidden "_GoBack" bookmark sitting
# right at the seam between the two pieces (this is where Word leaves the
# bookmark after the in-place deletion).

$prefix = " artisan serve "
$removed = "--host=0.0.0.0 "
$suffix = "--port="

# Step 1: remove "--host=0.0.0.0 " from the run, leaving " artisan serve --port="
$rng = $d.Content
$found = $rng.Find.Execute($removed, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the '--host=0.0.0.0 ' text to remove"
}
$rng.Delete()

# Step 2: find the resulting merged phrase so we know exactly where the
# seam between " artisan serve " and "--port=" now sits.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($prefix + $suffix, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the merged ' artisan serve --port=' text"
}
$seam = $rng2.Start + $prefix.Length

# Step 3: relocate the hidden "_GoBack" bookmark to that seam. Adding a
# bookmark with a name that already exists simply moves it, and because
# the bookmark is a zero-length anchor sitting between two otherwise
# identical runs, Word keeps them as two separate <w:r> elements split
# right at the bookmark.
$seamRange = $d.Range($seam, $seam)
$d.Bookmarks.Add("_GoBack", $seamRange)
